$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (metric_2D, PC)
$ws.Range("C2").Value = 0.0693
$ws.Range("D2").Value = 0.0594
$ws.Range("F2").Value = 0.099

# Row 3 (metric_breast_cancer, PC)
$ws.Range("C3").Value = 0.08790000000000001
$ws.Range("D3").Value = 0.0562
$ws.Range("F3").Value = 0.0527

# Row 4 (metric_load_iris, PC)
$ws.Range("C4").Value = 0.0333
$ws.Range("D4").Value = 0.0133
$ws.Range("F4").Value = 0.0533

# Row 5 (metric_load_wine, PC)
$ws.Range("C5").Value = 0.0618
$ws.Range("F5").Value = 0.0281

# Row 6 (indices_PC_LabelCorrection_before_fix_OCPC)
$ws.Range("G6").Value = 0.06307500000000001
$ws.Range("H6").Value = 0.000385801875
$ws.Range("I6").Value = 0.01964183990872546

# Row 8 (indices_PC_LabelCorrection_after_fix_OCPC)
$ws.Range("J8").Value = 0.0243
$ws.Range("K8").Value = 0.00060106
$ws.Range("L8").Value = 0.02451652503924649

# Row 9 (indices_CL_after_fix_OCPC)
$ws.Range("J9").Value = 0.03885
$ws.Range("K9").Value = 0.0011921425
$ws.Range("L9").Value = 0.0345274166424307

# Row 10 (metric_2D, LOF)
$ws.Range("C10").Value = 0.0743
$ws.Range("D10").Value = 0.0693
$ws.Range("F10").Value = 0.0891

# Row 11 (metric_breast_cancer, LOF)
$ws.Range("C11").Value = 0.0914
$ws.Range("D11").Value = 0.0369
$ws.Range("F11").Value = 0.0404

# Row 12 (metric_load_iris, LOF)
$ws.Range("C12").Value = 0.08
$ws.Range("D12").Value = 0.0533
$ws.Range("F12").Value = 0.08

# Row 13 (metric_load_wine, LOF)
$ws.Range("C13").Value = 0.1236
$ws.Range("D13").Value = 0.0506
$ws.Range("F13").Value = 0.0337

# Row 14 (indices_PC_LabelCorrection_before_fix_LOF)
$ws.Range("G14").Value = 0.092325
$ws.Range("H14").Value = 0.000363946875
$ws.Range("I14").Value = 0.01907739172423736

# Row 16 (indices_PC_LabelCorrection_after_fix_LOF)
$ws.Range("J16").Value = 0.03501666666666667
$ws.Range("K16").Value = 0.0007013913888888889
$ws.Range("L16").Value = 0.02648379483550061

# Row 17 (indices_CL_after_fix_LOF)
$ws.Range("J17").Value = 0.04053333333333334
$ws.Range("K17").Value = 0.001208158888888889
$ws.Range("L17").Value = 0.03475858007584443
